$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3 (Hydrogen / Non-metallic minerals) value removed -> blank cell
$ws.Range("D3").Value = ""

# C4 (Methanol / Chemicals) corrected value
$ws.Range("C4").Value = 46.31172096148298

# C5 (Ammonia / Chemicals) corrected value
$ws.Range("C5").Value = 3536.063524208197

# Row 7 renamed from "Other" to "Biogas", with a corrected D value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 421.3978040552909

# New row 8: "Other", copying row 7's formatting for column A, value in D
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 174.2199394676286
